$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# ---- Schedule sheet ----
$wsSchedule.Range("A2").Value = 46040.27083333334
$wsSchedule.Range("B2").Value = 46040.85416666666
$wsSchedule.Range("C2").Value = 14
$wsSchedule.Range("D2").Value = 52.91999999999999
$wsSchedule.Range("E2").Value = 118.342185
$wsSchedule.Range("F2").Value = 2.236246882086168
$wsSchedule.Range("A3").Value = 46040.91666666666
$wsSchedule.Range("C3").Value = 5
$wsSchedule.Range("D3").Value = 18.9
$wsSchedule.Range("E3").Value = 427.33517775
$wsSchedule.Range("F3").Value = 22.61032686507937
$wsSchedule.Range("E4").Value = 9.549998250000009
$wsSchedule.Range("F4").Value = 0.2807171737213407

# ---- Detailed sheet ----
$wsDetailed.Range("E14").Value = "OFF"
$wsDetailed.Range("B39").Value = 4.23686
$wsDetailed.Range("B40").Value = 56.85524
$wsDetailed.Range("B41").Value = 56.98
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("E41").Value = "ON"
$wsDetailed.Range("B42").Value = 56.98
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("E42").Value = "ON"
$wsDetailed.Range("B43").Value = 57.3
$wsDetailed.Range("B44").Value = 57.8068
$wsDetailed.Range("B45").Value = 57.3
$wsDetailed.Range("E45").Value = "OFF"
$wsDetailed.Range("B46").Value = 56.98
$wsDetailed.Range("B47").Value = 47.11292
$wsDetailed.Range("B48").Value = 36.2
$wsDetailed.Range("B49").Value = 36.2
$wsDetailed.Range("B50").Value = 36.2
$wsDetailed.Range("B51").Value = 36.2
$wsDetailed.Range("B52").Value = 36.2
$wsDetailed.Range("B54").Value = 47.79916
$wsDetailed.Range("B55").Value = 48.42041
$wsDetailed.Range("B56").Value = 49.90902
$wsDetailed.Range("B59").Value = 57.06003
$wsDetailed.Range("B60").Value = 58.29092
$wsDetailed.Range("B61").Value = 60.17508
$wsDetailed.Range("B65").Value = 36.06
$wsDetailed.Range("B66").Value = 0.7
$wsDetailed.Range("B68").Value = 0.01003
$wsDetailed.Range("B69").Value = -5.6862
$wsDetailed.Range("B70").Value = -6.14423
$wsDetailed.Range("B71").Value = -6.11853
$wsDetailed.Range("B74").Value = -5.50985
$wsDetailed.Range("B75").Value = -5.50985
$wsDetailed.Range("B76").Value = -6.01122
$wsDetailed.Range("B77").Value = -6.07128
$wsDetailed.Range("B78").Value = -5.01
$wsDetailed.Range("B79").Value = -4.70876
$wsDetailed.Range("B80").Value = -2.57526
$wsDetailed.Range("B82").Value = 0.0003
$wsDetailed.Range("B83").Value = -2.48809
$wsDetailed.Range("B84").Value = 0.00976
$wsDetailed.Range("B85").Value = 0.40914
$wsDetailed.Range("B86").Value = 12.60083
$wsDetailed.Range("B87").Value = 44.45378
$wsDetailed.Range("B88").Value = 57.38802
$wsDetailed.Range("B89").Value = 77.18000000000001
$wsDetailed.Range("B90").Value = 77.61297999999999
$wsDetailed.Range("B92").Value = 65
$wsDetailed.Range("B94").Value = 61.17796
$wsDetailed.Range("B95").Value = 58.65468
$wsDetailed.Range("B96").Value = 57.64573
$wsDetailed.Range("B97").Value = 61.65027
